$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6141.5713
$ws.Range("I18").Value = 3831.5
$ws.Range("K18").Value = 3831.5
$ws.Range("M18").Value = -3547.5
$ws.Range("H28").Value = 684.5714
$ws.Range("I28").Value = 382
$ws.Range("K28").Value = 382
$ws.Range("M28").Value = 103
$ws.Range("H40").Value = 1424.25
$ws.Range("J40").Value = 1424.25
$ws.Range("L40").Value = 1424.25
$ws.Range("N40").Value = -1774.25
$ws.Range("H62").Value = 3102.55
$ws.Range("I62").Value = 2783.875
$ws.Range("J62").Value = 3315
$ws.Range("K62").Value = 2783.875
$ws.Range("L62").Value = 3315
$ws.Range("M62").Value = -2159.875
$ws.Range("N62").Value = -4563
$ws.Range("H64").Value = 4319.1113
$ws.Range("H65").Value = 3102.55
$ws.Range("I65").Value = 2783.875
$ws.Range("J65").Value = 3315
$ws.Range("K65").Value = 13919.375
$ws.Range("L65").Value = 16575
$ws.Range("M65").Value = -10799.375
$ws.Range("N65").Value = -22815
$ws.Range("H67").Value = 4319.1113
$ws.Range("H88").Value = 1349.5
$ws.Range("J88").Value = 1382.6666
$ws.Range("L88").Value = 1382.6666
$ws.Range("N88").Value = -2194.6666
$ws.Range("H91").Value = 1349.5
$ws.Range("J91").Value = 1382.6666
$ws.Range("L91").Value = 1382.6666
$ws.Range("N91").Value = -4190.6666
$ws.Range("H137").Value = 2065.6
$ws.Range("I137").Value = 1899.1428
$ws.Range("J137").Value = 2211.25
$ws.Range("K137").Value = 5697.428400000001
$ws.Range("L137").Value = 6633.75
$ws.Range("M137").Value = -3147.428400000001
$ws.Range("N137").Value = -11733.75
$ws.Range("H138").Value = 1968.5
$ws.Range("I138").Value = 1748
$ws.Range("K138").Value = 5244
$ws.Range("M138").Value = -104
$ws.Range("H141").Value = 3094.724
$ws.Range("I141").Value = 2665.1924
$ws.Range("K141").Value = 7995.5772
$ws.Range("M141").Value = -2815.5772

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5870.778
$ws.Range("I74").Value = 5274.7144
$ws.Range("K74").Value = 5274.7144
$ws.Range("M74").Value = -4400.7144
$ws.Range("H77").Value = 5870.778
$ws.Range("I77").Value = 5274.7144
$ws.Range("K77").Value = 26373.572
$ws.Range("M77").Value = -22005.572
$ws.Range("H132").Value = 1916.6666
$ws.Range("J132").Value = 1950
$ws.Range("L132").Value = 5850
$ws.Range("N132").Value = -10910

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 187.85715
$ws.Range("I7").Value = 41.25
$ws.Range("J7").Value = 383.33334
$ws.Range("K7").Value = 41.25
$ws.Range("L7").Value = 383.33334
$ws.Range("M7").Value = 71.75
$ws.Range("N7").Value = -609.33334
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620
$ws.Range("H82").Value = 20016.5
$ws.Range("I82").Value = 6694.1665
$ws.Range("K82").Value = 6694.1665
$ws.Range("M82").Value = -6311.1665
$ws.Range("H85").Value = 20016.5
$ws.Range("I85").Value = 6694.1665
$ws.Range("K85").Value = 6694.1665
$ws.Range("M85").Value = -5368.1665
$ws.Range("H86").Value = 2815.8333
$ws.Range("I86").Value = 2347.25
$ws.Range("K86").Value = 2347.25
$ws.Range("M86").Value = -1224.25
$ws.Range("H89").Value = 2815.8333
$ws.Range("I89").Value = 2347.25
$ws.Range("K89").Value = 11736.25
$ws.Range("M89").Value = -6120.25
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H134").Value = 12834.2
$ws.Range("I134").Value = 11608.071
$ws.Range("K134").Value = 34824.213
$ws.Range("M134").Value = -32289.213

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4694.8125
$ws.Range("I31").Value = 3523.7
$ws.Range("K31").Value = 3523.7
$ws.Range("M31").Value = -3228.7
$ws.Range("H34").Value = 4694.8125
$ws.Range("I34").Value = 3523.7
$ws.Range("K34").Value = 3523.7
$ws.Range("M34").Value = -3321.7
$ws.Range("H68").Value = 29756.643
$ws.Range("H71").Value = 29756.643
$ws.Range("H99").Value = 3000
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H126").Value = 3000
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 7237.5
$ws.Range("J132").Value = 7000
$ws.Range("L132").Value = 21000
$ws.Range("N132").Value = -26060
$ws.Range("H134").Value = 3913.2666
$ws.Range("I134").Value = 3834.6428
$ws.Range("K134").Value = 11503.9284
$ws.Range("M134").Value = -8968.928400000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 916.6667
$ws.Range("H36").Value = 258.4
$ws.Range("I36").Value = 243
$ws.Range("K36").Value = 729
$ws.Range("M36").Value = -560
$ws.Range("H121").Value = 1054.6471
$ws.Range("I121").Value = 974.75
$ws.Range("K121").Value = 2924.25
$ws.Range("M121").Value = -1614.25
$ws.Range("H131").Value = 1725.091
$ws.Range("J131").Value = 1997.375
$ws.Range("L131").Value = 5992.125
$ws.Range("N131").Value = -16072.125
$ws.Range("H139").Value = 2254.111
$ws.Range("I139").Value = 2254.111
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 6762.333
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -1622.333
$ws.Range("N139").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 649.1539
$ws.Range("J97").Value = 1003
$ws.Range("L97").Value = 1003
$ws.Range("N97").Value = -1995
$ws.Range("H102").Value = 1197.7142
$ws.Range("I102").Value = 1197.7142
$ws.Range("K102").Value = 1197.7142
$ws.Range("M102").Value = 424.2858000000001
$ws.Range("H107").Value = 64.666664
$ws.Range("I107").Value = 64.666664
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 64.666664
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1855.333336
$ws.Range("N107").ClearContents()
$ws.Range("H126").Value = 1999.5
$ws.Range("I126").Value = 1999.5
$ws.Range("K126").Value = 5998.5
$ws.Range("M126").Value = -3528.5
$ws.Range("H132").Value = 3731.1
$ws.Range("I132").Value = 3473
$ws.Range("J132").Value = 4333.3335
$ws.Range("K132").Value = 10419
$ws.Range("L132").Value = 13000.0005
$ws.Range("M132").Value = -7889
$ws.Range("N132").Value = -18060.0005

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1990.05
$ws.Range("I22").Value = 1961.1666
$ws.Range("J22").Value = 2033.375
$ws.Range("K22").Value = 1961.1666
$ws.Range("L22").Value = 2033.375
$ws.Range("M22").Value = -1666.1666
$ws.Range("N22").Value = -2623.375
$ws.Range("H27").Value = 1990.05
$ws.Range("I27").Value = 1961.1666
$ws.Range("J27").Value = 2033.375
$ws.Range("K27").Value = 1961.1666
$ws.Range("L27").Value = 2033.375
$ws.Range("M27").Value = -1854.1666
$ws.Range("N27").Value = -2247.375
$ws.Range("H46").Value = 3570.9583
$ws.Range("I46").Value = 3299.923
$ws.Range("K46").Value = 3299.923
$ws.Range("M46").Value = -3111.923
$ws.Range("H62").Value = 49999.668
$ws.Range("J62").Value = 49999.668
$ws.Range("L62").Value = 49999.668
$ws.Range("N62").Value = -51247.668
$ws.Range("H65").Value = 49999.668
$ws.Range("J65").Value = 49999.668
$ws.Range("L65").Value = 149999.004
$ws.Range("N65").Value = -156239.004
$ws.Range("H122").Value = 3031.4
$ws.Range("I122").Value = 2851.75
$ws.Range("K122").Value = 8555.25
$ws.Range("M122").Value = -6105.25
$ws.Range("H132").Value = 10654.444
$ws.Range("I132").Value = 5200
$ws.Range("J132").Value = 17472.5
$ws.Range("K132").Value = 15600
$ws.Range("L132").Value = 52417.5
$ws.Range("M132").Value = -13070
$ws.Range("N132").Value = -57477.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1048.5
$ws.Range("I122").Value = 1048.5
$ws.Range("K122").Value = 3145.5
$ws.Range("M122").Value = -695.5
$ws.Range("H132").Value = 2416
$ws.Range("I132").Value = 2416
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7248
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4718
$ws.Range("N132").ClearContents()
